$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("main")
$ws2 = $wb.Worksheets.Item("search")
$ws3 = $wb.Worksheets.Item("week_week")

# --- Value updates -------------------------------------------------------

# main (sheet1): inital/final date text cells
$ws1.Range("D5").Value = "2015-02-07"
$ws1.Range("C5").Value = "2015-02-01"
$ws1.Range("D12").Value = $true

# search (sheet2): inital/final date text cells + table choice
$ws2.Range("C5").Value = "2015-01-01"
$ws2.Range("D5").Value = "2015-12-02"
$ws2.Range("G5").Value = "mFRR_Energy"

# week_week (sheet3): inital date text cell + number of weeks + flags
$ws3.Range("C5").Value = "2024-05-29"
$ws3.Range("D5").Value = 5
$ws3.Range("D10").Value = $true
$ws3.Range("E10").Value = $true
$ws3.Range("D12").Value = $false
$ws3.Range("E12").Value = $false

# --- Selection / active-cell updates -------------------------------------
# Touch sheets in this order so the last one (week_week) ends up active,
# matching the workbook's saved tabSelected/activeTab state.

$ws1.Range("E21").Select() | Out-Null
$ws2.Range("F5").Select() | Out-Null
$ws3.Range("F22").Select() | Out-Null
